# Revert capacity charts to show kilowatts on the y-axis.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the raw data values (column E, "Solar") for rows 21-23 from Watts to Kilowatts
$ws.Range("E21").Value = 11
$ws.Range("E22").Value = 9.300000000000001
$ws.Range("E23").Value = 16.82

# Update the number format for these cells (shared numFmt 164: "#,##0" -> "#,##0.0")
$ws.Range("E21:E23").NumberFormat = "#,##0.0"

# Update the chart: axis title text and number format
$chart = $wb.Worksheets.Item("Sheet1").ChartObjects(1).Chart
$valAxis = $chart.Axes(2)  # xlValue
$valAxis.AxisTitle.Text = "Kilowatts (kW)"
$valAxis.TickLabels.NumberFormat = "#,##0"
